# Add the "Jett:" / jett:cellRef(...) example row and the "Static:" /
# static-method-calling example row to the ExprTest worksheet, just below
# the existing "List:" example (row 15), matching the template update
# that accompanied the JEXL static-method-resolution feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: Jett custom function example
$ws.Range("A17").Value = "Jett:"
$ws.Range("B17").Value = "`${jett:cellRef(16, 1)}"
$ws.Range("D17").Value = "`${jett:cellRef(16, 1, 2, 3)}"

# Row 18: JEXL static method calling example
$ws.Range("A18").Value = "Static:"
$ws.Range("B18").Value = "`${java.lang.String.format('%s supports static method calling!', testBean2)}"
